$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = "img"
$ws.Range("O2").Value = "001.png"
$ws.Range("O3").Value = "002.png"

$ws.Range("O4").Select()
